# Apply stock-report corrections: several item rows get revised Qty (F) and
# Value (G) figures, a handful of rows have their Code/Rate/Qty/Value entries
# swapped with the following row, and all affected "Sub Total" / "Grand Total"
# rows are recalculated to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F70").Value = 2
$ws.Range("G70").Value = 269.9
$ws.Range("F71").Value = 297
$ws.Range("G71").Value = 18918.9
$ws.Range("F77").Value = 234
$ws.Range("G77").Value = 10937.16
$ws.Range("B90").Value = 165779.14
$ws.Range("F144").Value = 944
$ws.Range("G144").Value = 7976.8
$ws.Range("F145").Value = 382
$ws.Range("G145").Value = 3052.18
$ws.Range("B147").Value = 12207.64
$ws.Range("F151").Value = 87
$ws.Range("G151").Value = 7558.56
$ws.Range("B156").Value = 28537.83
$ws.Range("F229").Value = 51
$ws.Range("G229").Value = 7317.48
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("F249").Value = 134
$ws.Range("G249").Value = 18467.88
$ws.Range("B260").Value = 165784.64
$ws.Range("F296").Value = 22
$ws.Range("G296").Value = 466.4
$ws.Range("B304").Value = 161258.2
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B366").Value = 65066
$ws.Range("E366").Value = 13.61
$ws.Range("F366").Value = 90
$ws.Range("G366").Value = 1152.9
$ws.Range("B367").Value = 53263
$ws.Range("E367").Value = 15.29
$ws.Range("F367").Value = -309
$ws.Range("G367").Value = -3958.29
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("B382").Value = 45702
$ws.Range("E382").Value = 31.43
$ws.Range("F382").Value = -215
$ws.Range("G382").Value = -5654.5
$ws.Range("B383").Value = 64919
$ws.Range("E383").Value = 27.97
$ws.Range("F383").Value = 61
$ws.Range("G383").Value = 1604.3
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("B442").Value = 64810
$ws.Range("E442").Value = 291.22
$ws.Range("F442").Value = 4
$ws.Range("G442").Value = 1095.68
$ws.Range("B443").Value = 53319
$ws.Range("E443").Value = 310.64
$ws.Range("F443").Value = -6
$ws.Range("G443").Value = -1643.52
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 105
$ws.Range("G473").Value = 3447.15
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("B572").Value = 65079
$ws.Range("F572").Value = 6
$ws.Range("G572").Value = 245.22
$ws.Range("B573").Value = 65362
$ws.Range("F573").Value = 18
$ws.Range("G573").Value = 735.66
$ws.Range("F599").Value = 1288
$ws.Range("G599").Value = 210085.68
$ws.Range("F601").Value = 355
$ws.Range("G601").Value = 100418.85
$ws.Range("F602").Value = 306
$ws.Range("G602").Value = 44262.9
$ws.Range("B606").Value = 355615.48
$ws.Range("F610").Value = 10
$ws.Range("G610").Value = 409.9
$ws.Range("B618").Value = 41071.53
$ws.Range("B619").Value = 1553309.47
$ws.Range("B620").Value = 1553309.47
